$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, used to restore style
# after forcing NumberFormat to text ("@") on numeric-looking D-column values.
$plainStyle = $ws.Range("B2").Style

# --- Column D (Price) updates ---
# NumberFormat is forced to text first so Excel does not silently
# reinterpret these dotted/decimal-looking strings as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.025.43"
$ws.Range("D2").Style = $plainStyle
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.051.95"
$ws.Range("D3").Style = $plainStyle
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.52"
$ws.Range("D5").Style = $plainStyle
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.59"
$ws.Range("D7").Style = $plainStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.20"
$ws.Range("D9").Style = $plainStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").Style = $plainStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.43"
$ws.Range("D13").Style = $plainStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.353.00"
$ws.Range("D14").Style = $plainStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.805"
$ws.Range("D15").Style = $plainStyle
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.056.49"
$ws.Range("D17").Style = $plainStyle
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.970.35"
$ws.Range("D18").Style = $plainStyle
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.67"
$ws.Range("D19").Style = $plainStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.36"
$ws.Range("D20").Style = $plainStyle
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0911"
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.45"
$ws.Range("D22").Style = $plainStyle
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.97"
$ws.Range("D23").Style = $plainStyle
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.24"
$ws.Range("D27").Style = $plainStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("D28").Style = $plainStyle
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("D31").Style = $plainStyle
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.74"
$ws.Range("D32").Style = $plainStyle
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").Style = $plainStyle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.26"
$ws.Range("D37").Style = $plainStyle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.114"
$ws.Range("D39").Style = $plainStyle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.83"
$ws.Range("D41").Style = $plainStyle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0224"
$ws.Range("D42").Style = $plainStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.35"
$ws.Range("D44").Style = $plainStyle
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.87"
$ws.Range("D45").Style = $plainStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.58"
$ws.Range("D46").Style = $plainStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.48"
$ws.Range("D47").Style = $plainStyle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.290.89"
$ws.Range("D48").Style = $plainStyle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.87"
$ws.Range("D50").Style = $plainStyle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.239.86"
$ws.Range("D51").Style = $plainStyle

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +8.31%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +5.23%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  +8.46%  "
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("E16").Value = "  +7.13%  "
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  +14.90%  "
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("E21").Value = "  +6.57%  "
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +12.05%  "
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -3.27%  "
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("E39").Value = "  +17.81%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("E46").Value = "  +12.73%  "
$ws.Range("E47").Value = "  +4.91%  "
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("E51").Value = "  -2.16%  "
